# Apply the changes described by the commit: selection/active-cell changes
# on "TakeAction" and "Otcandnonrx" sheets, and content changes to A5/B5 on
# "Otcandnonrx".

$wb = $excel.ActiveWorkbook

# --- Sheet "TakeAction": move selection from A3:B6 to A12 ---
$wsTakeAction = $wb.Worksheets.Item("TakeAction")
$wsTakeAction.Activate()
$wsTakeAction.Range("A12").Select()

# --- Sheet "Otcandnonrx": update A5/B5 values and move selection ---
$wsOtc = $wb.Worksheets.Item("Otcandnonrx")
$wsOtc.Activate()

# Update the cell values (also clears the previous cell style, matching
# the target workbook which drops the explicit style index on A5/B5).
$wsOtc.Range("A5").ClearFormats()
$wsOtc.Range("B5").ClearFormats()
$wsOtc.Range("A5").Value = "Cura Tulsi Ark Drops 50 ml"
$wsOtc.Range("B5").Value = "Cura Tulsi Ark Drops 50 ml"

# Move the active selection to B11
$wsOtc.Range("B11").Select()

# Re-activate Otcandnonrx sheet so it remains the tab shown (tabSelected)
$wsOtc.Activate()
